$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (e.g. "1.001", "29.809.71") that must
# stay plain text. Pre-format the affected D cells as Text so assigning the
# new strings does not get reinterpreted as a number, then restore the
# "Normal" style so the on-disk style index matches the original (style 0).
$dCells = @("D2","D3","D4","D6","D7","D8","D9","D10","D11","D12","D14","D15","D16","D17","D18","D20","D21","D22","D23","D24","D26","D27","D29","D30","D31","D34","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D47","D49","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "29.809.71"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "1.872.19"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("E5").Value = "  -1.25%  "
$ws.Range("D6").Value = "241.31"
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "0.3138"
$ws.Range("E8").Value = "  -0.56%  "
$ws.Range("D9").Value = "0.07128"
$ws.Range("E9").Value = "  -0.59%  "
$ws.Range("D10").Value = "24.41"
$ws.Range("E10").Value = "  -1.66%  "
$ws.Range("D11").Value = "0.08154"
$ws.Range("E11").Value = "  -3.28%  "
$ws.Range("D12").Value = "1.895.17"
$ws.Range("E12").Value = "  +1.51%  "
$ws.Range("E13").Value = "  -1.82%  "
$ws.Range("D14").Value = "5.341"
$ws.Range("E14").Value = "  -1.40%  "
$ws.Range("D15").Value = "92.45"
$ws.Range("E15").Value = "  -0.44%  "
$ws.Range("D16").Value = "29.847.91"
$ws.Range("E16").Value = "  -0.40%  "
$ws.Range("D17").Value = "6.003"
$ws.Range("E17").Value = "  -1.63%  "
$ws.Range("D18").Value = "248.35"
$ws.Range("E18").Value = "  +1.86%  "
$ws.Range("E19").Value = "  -1.94%  "
$ws.Range("D20").Value = "0.000007806"
$ws.Range("E20").Value = "  -0.35%  "
$ws.Range("D21").Value = "2.159.64"
$ws.Range("E21").Value = "  +2.10%  "
$ws.Range("D22").Value = "1.003"
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  -0.28%  "
$ws.Range("D24").Value = "7.751"
$ws.Range("E24").Value = "  -3.06%  "
$ws.Range("E25").Value = "  -1.54%  "
$ws.Range("D26").Value = "9.208"
$ws.Range("E26").Value = "  -1.25%  "
$ws.Range("D27").Value = "164.18"
$ws.Range("E27").Value = "  -0.36%  "
$ws.Range("E28").Value = "  -0.81%  "
$ws.Range("D29").Value = "2.017"
$ws.Range("E29").Value = "  -1.30%  "
$ws.Range("D30").Value = "1.449"
$ws.Range("D31").Value = "4.522"
$ws.Range("E31").Value = "  -1.99%  "
$ws.Range("E32").Value = "  -0.73%  "
$ws.Range("E33").Value = "  -2.84%  "
$ws.Range("D34").Value = "0.05311"
$ws.Range("E34").Value = "  -0.58%  "
$ws.Range("E35").Value = "  -0.86%  "
$ws.Range("D36").Value = "0.7406"
$ws.Range("E36").Value = "  -2.61%  "
$ws.Range("D37").Value = "1.001"
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("D38").Value = "2.705"
$ws.Range("E38").Value = "  +0.21%  "
$ws.Range("D39").Value = "0.01938"
$ws.Range("E39").Value = "  -0.98%  "
$ws.Range("D40").Value = "2.732"
$ws.Range("E40").Value = "  -0.74%  "
$ws.Range("D41").Value = "0.4465"
$ws.Range("E41").Value = "  -0.64%  "
$ws.Range("D42").Value = "5.964"
$ws.Range("E42").Value = "  -3.04%  "
$ws.Range("D43").Value = "0.8690"
$ws.Range("E43").Value = "  +0.56%  "
$ws.Range("D44").Value = "71.23"
$ws.Range("E44").Value = "  -2.09%  "
$ws.Range("D45").Value = "1.047.61"
$ws.Range("E45").Value = "  -6.12%  "
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("D47").Value = "103.92"
$ws.Range("E47").Value = "  +0.61%  "
$ws.Range("E48").Value = "  -1.64%  "
$ws.Range("D49").Value = "7.427"
$ws.Range("E49").Value = "  -3.68%  "
$ws.Range("D50").Value = "2.054.49"
$ws.Range("D51").Value = "9.507"
$ws.Range("E51").Value = "  -0.49%  "

foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}
